$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (D) text or $null to leave unchanged, new Volume(1h) (E) text
$updates = @(
    ,@(2, "37.160.42", "  -0.35%  ")
    ,@(3, "2.074.46", "  -1.04%  ")
    ,@(4, $null, "  -0.04%  ")
    ,@(5, "253.59", "  +1.18%  ")
    ,@(6, $null, "  +1.91%  ")
    ,@(7, "58.99", "  +9.84%  ")
    ,@(8, $null, "  +0.01%  ")
    ,@(9, $null, "  +4.34%  ")
    ,@(10, "61.43", "  -0.55%  ")
    ,@(11, "0.0799", "  +7.68%  ")
    ,@(12, $null, "  +2.67%  ")
    ,@(13, $null, "  +6.83%  ")
    ,@(14, "2.379.29", "  -0.97%  ")
    ,@(15, "0.823", "  -1.91%  ")
    ,@(16, $null, "  +7.47%  ")
    ,@(17, "2.075.15", "  -0.95%  ")
    ,@(18, "37.169.20", "  -0.27%  ")
    ,@(19, "15.79", "  +8.44%  ")
    ,@(20, "74.73", "  +2.90%  ")
    ,@(21, "0.0₃0927", "  +10.49%  ")
    ,@(22, $null, "  +5.07%  ")
    ,@(23, "240.31", "  -0.33%  ")
    ,@(24, $null, "  +0.00%  ")
    ,@(25, $null, "  -3.04%  ")
    ,@(26, "2.29", "  +13.53%  ")
    ,@(27, "169.97", "  -1.02%  ")
    ,@(28, "9.35", "  +0.51%  ")
    ,@(29, "20.36", "  -1.52%  ")
    ,@(30, "0.126", "  +2.86%  ")
    ,@(31, $null, "  +6.53%  ")
    ,@(32, "4.80", "  +6.35%  ")
    ,@(33, "0.0634", "  +3.37%  ")
    ,@(34, "4.50", "  +9.38%  ")
    ,@(35, "0.0907", "  +0.06%  ")
    ,@(36, $null, "  -0.05%  ")
    ,@(37, "2.30", "  -0.11%  ")
    ,@(38, $null, "  -3.83%  ")
    ,@(39, $null, "  +26.06%  ")
    ,@(40, $null, "  +1.90%  ")
    ,@(41, "0.0228", "  +1.08%  ")
    ,@(42, "17.84", "  -3.80%  ")
    ,@(43, $null, "  +0.45%  ")
    ,@(44, "99.15", "  +0.31%  ")
    ,@(45, "4.43", "  +7.49%  ")
    ,@(46, $null, "  +1.01%  ")
    ,@(47, "4.58", "  +14.27%  ")
    ,@(48, $null, "  +8.64%  ")
    ,@(49, "1.304.77", "  -1.01%  ")
    ,@(50, "2.93", "  -1.03%  ")
    ,@(51, $null, "  -1.55%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]
    if ($null -ne $dVal) {
        # Force text so numeric-looking strings (e.g. "4.80") keep their
        # exact text, matching the original inlineStr cell, then restore
        # the default "Normal" style so no stray cell format is introduced.
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $dVal
        $cell.Style = "Normal"
    }
    if ($null -ne $eVal) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $eVal
        $cell.Style = "Normal"
    }
}

Write-Output "Updated $($updates.Count) rows"
